# "Generate Report for Handoff" - refresh the localization-status report:
#  - flip status from "Handed back: in sync with en-US" to "Ready for handoff"
#    on the Overview sheet (both language columns) and on each language
#    sheet's Status column
#  - bump the handoff timestamps to the new generation time
#  - the Status column is narrower now that the text is shorter (re-fit)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# Overview sheet: E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-31 13:03:43"

# zh-cn sheet: C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-31 13:03:31"

# de-de sheet: C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-31 13:03:43"

# Re-fit the now-shorter Status columns
$wsOverview.Columns.Item(5).ColumnWidth = 16.45
$wsOverview.Columns.Item(6).ColumnWidth = 16.45
$wsZhCn.Columns.Item(3).ColumnWidth = 16.45
$wsDeDe.Columns.Item(3).ColumnWidth = 16.45
